$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the scattered cells that were removed from the sheet.
# (Row 5 had no value to begin with, but is included for completeness.)
$ws.Range("A1:A3").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("A12:A13").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A21:A24").ClearContents()
$ws.Range("A26:A27").ClearContents()
$ws.Range("A30").ClearContents()
$ws.Range("A35:A36").ClearContents()
$ws.Range("A39:A41").ClearContents()
$ws.Range("A225").ClearContents()

# Row 29 ends up as a bare, empty row (it sits between the cleared A26:A27 and
# A30 cells). Touching its outline level keeps an explicit <row r="29"/>
# element in the saved worksheet, matching the target structure.
$ws.Rows("29").OutlineLevel = 0

# Update the view: the window had scrolled back to the top and the selected
# cell moved to E9 (no more frozen/offset top-left cell at A252 / D258).
$ws.Range("E9").Select()
